# Auto-generated edit script: update column F ('想去人数') values
# per the commit diff (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$exhibitions = @{
    2 = 1377
    3 = 2231
    6 = 71
    7 = 698
    8 = 125
    11 = 2576
    12 = 1623
    13 = 1613
    15 = 266
    16 = 648
    17 = 824
    18 = 104
    19 = 332
    20 = 1094
    22 = 38
    23 = 533
    24 = 5588
    25 = 230
    26 = 903
    27 = 109
    29 = 146
    30 = 248
    32 = 46
    33 = 1061
    34 = 799
    36 = 64
    38 = 419
    39 = 1160
    40 = 146
    41 = 114
    42 = 191
    44 = 103
}
foreach ($row in $exhibitions.Keys) {
    $ws.Cells.Item($row, 6).Value = $exhibitions[$row]
}

$ws = $wb.Worksheets.Item("演出")
$shows = @{
    3 = 793
    5 = 434
    10 = 4
}
foreach ($row in $shows.Keys) {
    $ws.Cells.Item($row, 6).Value = $shows[$row]
}

$ws = $wb.Worksheets.Item("全部类型")
$allTypes = @{
    2 = 1377
    4 = 2231
    8 = 71
    9 = 698
    10 = 125
    15 = 2576
    16 = 1623
    17 = 1613
    19 = 266
    20 = 648
    22 = 824
    23 = 104
    24 = 332
    25 = 1094
    26 = 38
    27 = 533
    28 = 5588
    29 = 230
    30 = 904
    31 = 109
    33 = 146
    34 = 248
    36 = 46
    37 = 1061
    38 = 799
    39 = 64
    40 = 419
    41 = 1160
    42 = 146
    43 = 114
    44 = 191
    46 = 103
    47 = 4
}
foreach ($row in $allTypes.Keys) {
    $ws.Cells.Item($row, 6).Value = $allTypes[$row]
}
